$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 607 (the "コーヒーは急には飲まれない" quote) and shift all subsequent rows up.
$ws.Rows.Item(607).Delete() | Out-Null
